$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 3.1
$ws.Range("H2").Value = 2.44
$ws.Range("I2").Value = 2.74
$ws.Range("J2").Value = 3.15
$ws.Range("O2").Value = 1.44
$ws.Range("S2").Value = 3.85
$ws.Range("V2").Value = 1.57
$ws.Range("X2").Value = 13.5
$ws.Range("Y2").Value = 10.5
$ws.Range("AB2").Value = 13.5
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 14
$ws.Range("AG2").Value = 18

# Row 3 updates
$ws.Range("F3").Value = 1.65
$ws.Range("U3").Value = 1.78

# Row 4 updates
$ws.Range("G4").Value = 1.7
$ws.Range("J4").Value = 1.01
$ws.Range("W4").Value = 2.42
